$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sort each quadrant block by FDR ascending (Table13, Table134, Table135) ---
# Table1 (A1:C2) has a single data row, no sort needed.

# Table13: header A5:C5, data A6:C13
$ws.Range("A6:C13").Sort($ws.Range("B6:B13"))

# Table134: header A15:C15, data A16:C19
$ws.Range("A16:C19").Sort($ws.Range("B16:B19"))

# Table135: header A22:C22, data A23:C28
$ws.Range("A23:C28").Sort($ws.Range("B23:B28"))

# --- Remove the extra blank separator rows so each block is separated by one blank row ---
# Before: blank rows 3,4 between Table1 and Table13 -> keep only one (delete row 4)
$ws.Rows("4").Delete()
# After the above delete, the second extra gap (originally rows 20,21) is now rows 19,20 -> delete row 20
$ws.Rows("20").Delete()

# --- Rename the table objects to match the new numbering ---
# (ListObjects collection order isn't stable across renames, so look each one
#  up by its current name rather than by a positional index.)
function Get-ListObjectByName($sheet, $name) {
  for ($i = 1; $i -le $sheet.ListObjects.Count; $i++) {
    $lo = $sheet.ListObjects.Item($i)
    if ($lo.Name -eq $name) {
      return $lo
    }
  }
  return $null
}

(Get-ListObjectByName $ws "Table134").Name = "Table14"
(Get-ListObjectByName $ws "Table135").Name = "Table15"

# --- Rename the "Gene set" column header to "gene_set" everywhere (shared header text) ---
$ws.Range("A1").Value = "gene_set"
$ws.Range("A4").Value = "gene_set"
$ws.Range("A14").Value = "gene_set"
$ws.Range("A20").Value = "gene_set"

# --- Fix the capitalization of the tRNA PROCESSING gene set label ---
$ws.Range("A7").Value = "tRNA PROCESSING"

# --- Re-apply the centered formatting to the blank separator rows between blocks ---
# (copy the header row's format only, so no stray/unused cell style gets created)
$ws.Range("A1:C1").Copy()
foreach ($r in @(3, 13, 19)) {
  $ws.Range("A" + $r + ":C" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# --- Update selection to match the recorded cursor position ---
$ws.Range("D10").Select()
